# Update "F" column (想去人数 / number of attendees) values across the four
# worksheets to match the regenerated data output.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1745
$ws.Cells.Item(3, 6).Value = 10190
$ws.Cells.Item(8, 6).Value = 1632
$ws.Cells.Item(9, 6).Value = 183
$ws.Cells.Item(10, 6).Value = 389
$ws.Cells.Item(12, 6).Value = 213
$ws.Cells.Item(14, 6).Value = 484
$ws.Cells.Item(15, 6).Value = 1187
$ws.Cells.Item(18, 6).Value = 14
$ws.Cells.Item(19, 6).Value = 95
$ws.Cells.Item(20, 6).Value = 360
$ws.Cells.Item(22, 6).Value = 330
$ws.Cells.Item(23, 6).Value = 104
$ws.Cells.Item(24, 6).Value = 1162
$ws.Cells.Item(27, 6).Value = 40
$ws.Cells.Item(28, 6).Value = 435
$ws.Cells.Item(29, 6).Value = 236
$ws.Cells.Item(31, 6).Value = 468
$ws.Cells.Item(32, 6).Value = 225
$ws.Cells.Item(33, 6).Value = 374
$ws.Cells.Item(34, 6).Value = 527
$ws.Cells.Item(35, 6).Value = 649
$ws.Cells.Item(36, 6).Value = 755
$ws.Cells.Item(37, 6).Value = 528
$ws.Cells.Item(38, 6).Value = 1271
$ws.Cells.Item(39, 6).Value = 815
$ws.Cells.Item(40, 6).Value = 385
$ws.Cells.Item(41, 6).Value = 337
$ws.Cells.Item(42, 6).Value = 19
$ws.Cells.Item(43, 6).Value = 357
$ws.Cells.Item(44, 6).Value = 80
$ws.Cells.Item(45, 6).Value = 354
$ws.Cells.Item(46, 6).Value = 80

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 37
$ws.Cells.Item(14, 6).Value = 103
$ws.Cells.Item(15, 6).Value = 65
$ws.Cells.Item(18, 6).Value = 1102
$ws.Cells.Item(20, 6).Value = 711
$ws.Cells.Item(22, 6).Value = 328
$ws.Cells.Item(24, 6).Value = 77
$ws.Cells.Item(31, 6).Value = 209
$ws.Cells.Item(35, 6).Value = 189
$ws.Cells.Item(43, 6).Value = 70

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 808
$ws.Cells.Item(6, 6).Value = 2529
$ws.Cells.Item(7, 6).Value = 4107
$ws.Cells.Item(10, 6).Value = 312
$ws.Cells.Item(11, 6).Value = 197

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1745
$ws.Cells.Item(3, 6).Value = 808
$ws.Cells.Item(4, 6).Value = 10190
$ws.Cells.Item(7, 6).Value = 4107
$ws.Cells.Item(9, 6).Value = 312
$ws.Cells.Item(10, 6).Value = 312
$ws.Cells.Item(12, 6).Value = 1632
$ws.Cells.Item(13, 6).Value = 183
$ws.Cells.Item(14, 6).Value = 389
$ws.Cells.Item(16, 6).Value = 213
$ws.Cells.Item(19, 6).Value = 1187
$ws.Cells.Item(23, 6).Value = 103
$ws.Cells.Item(24, 6).Value = 65
$ws.Cells.Item(25, 6).Value = 95
$ws.Cells.Item(26, 6).Value = 1102
$ws.Cells.Item(27, 6).Value = 360
$ws.Cells.Item(28, 6).Value = 330
$ws.Cells.Item(30, 6).Value = 1162
$ws.Cells.Item(32, 6).Value = 77
$ws.Cells.Item(34, 6).Value = 236
$ws.Cells.Item(35, 6).Value = 366
$ws.Cells.Item(36, 6).Value = 468
$ws.Cells.Item(38, 6).Value = 374
$ws.Cells.Item(39, 6).Value = 527
$ws.Cells.Item(40, 6).Value = 649
$ws.Cells.Item(41, 6).Value = 209
$ws.Cells.Item(42, 6).Value = 755
$ws.Cells.Item(43, 6).Value = 528
$ws.Cells.Item(44, 6).Value = 815
$ws.Cells.Item(45, 6).Value = 385
$ws.Cells.Item(46, 6).Value = 337
$ws.Cells.Item(48, 6).Value = 357
$ws.Cells.Item(49, 6).Value = 354
$ws.Cells.Item(50, 6).Value = 70
